# Auto-generated edit script: applies numeric corrections to the
# currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) on
# several rows across the ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets, per the
# scheduled-runner refresh described in the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H33").Value = 231
$ws.Range("I33").Value = 192
$ws.Range("J33").Value = 413
$ws.Range("K33").Value = 192
$ws.Range("L33").Value = 413
$ws.Range("M33").Value = 37
$ws.Range("N33").Value = -871

$ws.Range("H39").Value = 409.36365
$ws.Range("I39").Value = 33.77778
$ws.Range("J39").Value = 2099.5
$ws.Range("K39").Value = 101.33334
$ws.Range("L39").Value = 6298.5
$ws.Range("M39").Value = 194.66666
$ws.Range("N39").Value = -6890.5

$ws.Range("H137").Value = 428944.06
$ws.Range("I137").Value = 2631.5
$ws.Range("J137").Value = 560117.1
$ws.Range("K137").Value = 7894.5
$ws.Range("L137").Value = 1680351.3
$ws.Range("M137").Value = -5344.5
$ws.Range("N137").Value = -1685451.3

$ws.Range("H138").Value = 3916.6
$ws.Range("J138").Value = 6666.6665
$ws.Range("L138").Value = 19999.9995
$ws.Range("N138").Value = -30279.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()

$ws.Range("H8").Value = 1701.6666
$ws.Range("I8").Value = 1701.6666
$ws.Range("K8").Value = 1701.6666
$ws.Range("M8").Value = -1557.6666

$ws.Range("H11").Value = 5004375
$ws.Range("I11").Value = 10002000
$ws.Range("J11").Value = 6750
$ws.Range("K11").Value = 10002000
$ws.Range("L11").Value = 6750
$ws.Range("M11").Value = -10001856
$ws.Range("N11").Value = -7038

$ws.Range("H32").Value = 3144.1528
$ws.Range("I32").Value = 2325.4626
$ws.Range("K32").Value = 2325.4626
$ws.Range("M32").Value = -2038.4626

$ws.Range("H132").Value = 2124.9355
$ws.Range("I132").Value = 2133.724
$ws.Range("J132").Value = 1997.5
$ws.Range("K132").Value = 6401.172
$ws.Range("L132").Value = 5992.5
$ws.Range("M132").Value = -3871.172
$ws.Range("N132").Value = -11052.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1292.4375
$ws.Range("I16").Value = 723.25
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 723.25
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -436.25
$ws.Range("N16").Value = -3574

$ws.Range("H31").Value = 2018.3036
$ws.Range("I31").Value = 1406.4286
$ws.Range("J31").Value = 2222.262
$ws.Range("K31").Value = 1406.4286
$ws.Range("L31").Value = 2222.262
$ws.Range("M31").Value = -1111.4286
$ws.Range("N31").Value = -2812.262

$ws.Range("H34").Value = 2018.3036
$ws.Range("I34").Value = 1406.4286
$ws.Range("J34").Value = 2222.262
$ws.Range("K34").Value = 1406.4286
$ws.Range("L34").Value = 2222.262
$ws.Range("M34").Value = -1204.4286
$ws.Range("N34").Value = -2626.262

$ws.Range("H105").Value = 2481.1667
$ws.Range("I105").Value = 1541.3125
$ws.Range("K105").Value = 1541.3125
$ws.Range("M105").Value = 205.6875

$ws.Range("H113").Value = 1292.4375
$ws.Range("I113").Value = 723.25
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 723.25
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1446.75
$ws.Range("N113").Value = -7340

$ws.Range("H124").Value = 50000
$ws.Range("J124").Value = 50000
$ws.Range("L124").Value = 50000
$ws.Range("N124").Value = -54910

$ws.Range("H132").Value = 384089.25
$ws.Range("I132").Value = 1694.1936
$ws.Range("J132").Value = 4335504.5
$ws.Range("K132").Value = 5082.5808
$ws.Range("L132").Value = 13006513.5
$ws.Range("M132").Value = -2552.5808
$ws.Range("N132").Value = -13011573.5

$ws.Range("H134").Value = 28291.564
$ws.Range("I134").Value = 2690.889
$ws.Range("J134").Value = 335499.66
$ws.Range("K134").Value = 8072.667
$ws.Range("L134").Value = 1006498.98
$ws.Range("M134").Value = -5537.667
$ws.Range("N134").Value = -1011568.98

$ws.Range("H141").Value = 100000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 100000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3763
$ws.Range("J68").Value = 4147.5
$ws.Range("L68").Value = 12442.5
$ws.Range("N68").Value = -14064.5

$ws.Range("H71").Value = 3763
$ws.Range("J71").Value = 4147.5
$ws.Range("L71").Value = 37327.5
$ws.Range("N71").Value = -45439.5

$ws.Range("H107").Value = 1814.5714
$ws.Range("J107").Value = 2641.5
$ws.Range("L107").Value = 7924.5
$ws.Range("N107").Value = -11764.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 928.73334
$ws.Range("I107").Value = 1062.7273
$ws.Range("J107").Value = 560.25
$ws.Range("K107").Value = 1062.7273
$ws.Range("L107").Value = 560.25
$ws.Range("M107").Value = 857.2727
$ws.Range("N107").Value = -4400.25

$ws.Range("H121").Value = 45000
$ws.Range("J121").Value = 45000
$ws.Range("L121").Value = 45000
$ws.Range("N121").Value = -48494

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2176.9
$ws.Range("I68").Value = 1971
$ws.Range("J68").Value = 3000.5
$ws.Range("K68").Value = 1971
$ws.Range("L68").Value = 3000.5
$ws.Range("M68").Value = -1222
$ws.Range("N68").Value = -4498.5

$ws.Range("H71").Value = 2176.9
$ws.Range("I71").Value = 1971
$ws.Range("J71").Value = 3000.5
$ws.Range("K71").Value = 9855
$ws.Range("L71").Value = 15002.5
$ws.Range("M71").Value = -6111
$ws.Range("N71").Value = -22490.5

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H132").Value = 3452.9443
$ws.Range("I132").Value = 3452.9443
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10358.8329
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7828.832900000001
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()

$ws.Range("H113").Value = 1057.381
$ws.Range("I113").Value = 1320.2307
$ws.Range("J113").Value = 630.25
$ws.Range("K113").Value = 3960.6921
$ws.Range("L113").Value = 1890.75
$ws.Range("M113").Value = -1790.6921
$ws.Range("N113").Value = -6230.75
